# Automatic update of files.
# 1) Column C ("Förändrad") for every data row (2-33) moves from 2026-02-13
#    (serial 46066) to 2026-02-17 (serial 46070).
# 2) Rows 11-18 are re-sorted: the block is rotated up by two rows, i.e. the
#    data that used to live in rows 13-18 moves up to rows 11-16, and the
#    data that used to live in rows 11-12 wraps around to rows 17-18
#    (carrying the "Markägare" value in column F along with its row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Step 1: bump "Förändrad" (column C) on every data row from 46066 to 46070.
# ---------------------------------------------------------------------
for ($row = 2; $row -le 33; $row++) {
    $ws.Cells.Item($row, 3).Value = 46070
}

# ---------------------------------------------------------------------
# Step 2: capture rows 11-18 (columns A, B, F, G) before touching anything,
# then write back the rotated order.
# ---------------------------------------------------------------------
$colA = @{}
$colB = @{}
$colF = @{}
$colG = @{}

for ($row = 11; $row -le 18; $row++) {
    $colA[$row] = $ws.Cells.Item($row, 1).Value2
    $colB[$row] = $ws.Cells.Item($row, 2).Value2
    $colF[$row] = $ws.Cells.Item($row, 6).Value2
    $colG[$row] = $ws.Cells.Item($row, 7).Value2
}

# new row r (11..16) <- old row r+2 ; new row 17 <- old row 11 ; new row 18 <- old row 12
$order = @{
    11 = 13
    12 = 14
    13 = 15
    14 = 16
    15 = 17
    16 = 18
    17 = 11
    18 = 12
}

foreach ($destRow in $order.Keys) {
    $srcRow = $order[$destRow]

    $ws.Cells.Item($destRow, 1).Value = $colA[$srcRow]
    $ws.Cells.Item($destRow, 2).Value = $colB[$srcRow]

    $srcF = $colF[$srcRow]
    if ($srcF -eq $null -or $srcF -eq "") {
        $ws.Cells.Item($destRow, 6).Value = $null
    } else {
        $ws.Cells.Item($destRow, 6).Value = $srcF
    }

    $ws.Cells.Item($destRow, 7).Value = $colG[$srcRow]
}
